$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores plain-text numbers (e.g. "595.70", "0.0388").
# Assigning such a numeric-looking string straight to Range.Value makes Excel
# auto-convert the cell to a real Number, which silently drops significant
# trailing zeros and introduces floating-point artifacts. Force those cells to
# Text format first so the literal digit string from the update is preserved.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "63.002.85"
$ws.Range("E2").Value = "  -2.46%  "
$ws.Range("D3").Value = "3.124.85"
$ws.Range("E3").Value = "  -0.84%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "595.70"
$ws.Range("E5").Value = "  -2.08%  "
$ws.Range("D6").Value = "136.32"
$ws.Range("E6").Value = "  -5.39%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "3.115.55"
$ws.Range("E8").Value = "  -1.08%  "
$ws.Range("D9").Value = "0.515"
$ws.Range("E9").Value = "  -1.88%  "
$ws.Range("E10").Value = "  -3.92%  "
$ws.Range("D11").Value = "5.26"
$ws.Range("E11").Value = "  -3.10%  "
$ws.Range("E12").Value = "  -2.78%  "
$ws.Range("E13").Value = "  -4.96%  "
$ws.Range("D14").Value = "34.21"
$ws.Range("E14").Value = "  -3.46%  "
$ws.Range("D15").Value = "3.636.77"
$ws.Range("E15").Value = "  -0.95%  "
$ws.Range("E16").Value = "  +2.29%  "
$ws.Range("D17").Value = "63.051.12"
$ws.Range("E17").Value = "  -2.12%  "
$ws.Range("D18").Value = "3.119.17"
$ws.Range("E18").Value = "  -1.13%  "
$ws.Range("E19").Value = "  -2.97%  "
$ws.Range("D20").Value = "473.80"
$ws.Range("E20").Value = "  -1.50%  "
$ws.Range("E21").Value = "  -3.69%  "
$ws.Range("D22").Value = "0.697"
$ws.Range("E22").Value = "  -2.75%  "
$ws.Range("D23").Value = "7.67"
$ws.Range("E23").Value = "  -0.49%  "
$ws.Range("D24").Value = "86.01"
$ws.Range("E24").Value = "  +0.86%  "
$ws.Range("D25").Value = "12.93"
$ws.Range("E25").Value = "  -3.87%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  -1.67%  "
$ws.Range("D28").Value = "7.91"
$ws.Range("E28").Value = "  -6.33%  "
$ws.Range("D29").Value = "6.94"
$ws.Range("E29").Value = "  -3.50%  "
$ws.Range("E30").Value = "  +1.66%  "
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("D32").Value = "26.66"
$ws.Range("E32").Value = "  -1.35%  "
$ws.Range("D33").Value = "0.108"
$ws.Range("E33").Value = "  -5.99%  "
$ws.Range("E34").Value = "  -4.83%  "
$ws.Range("E35").Value = "  -3.09%  "
$ws.Range("E36").Value = "  -3.27%  "
$ws.Range("D37").Value = "51.95"
$ws.Range("E37").Value = "  -0.95%  "
$ws.Range("D38").Value = "0.0₃0700"
$ws.Range("E38").Value = "  -8.98%  "
$ws.Range("D39").Value = "0.0388"
$ws.Range("E39").Value = "  -1.96%  "
$ws.Range("D40").Value = "420.66"
$ws.Range("E40").Value = "  -5.98%  "
$ws.Range("E41").Value = "  -0.49%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "2.896.52"
$ws.Range("E42").Value = "  +0.85%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "2.68"
$ws.Range("E43").Value = "  -11.68%  "
$ws.Range("D44").Value = "0.114"
$ws.Range("E44").Value = "  -5.38%  "
$ws.Range("E45").Value = "  +1.48%  "
$ws.Range("D47").Value = "2.11"
$ws.Range("E47").Value = "  -5.54%  "
$ws.Range("D48").Value = "25.50"
$ws.Range("E48").Value = "  -2.83%  "
$ws.Range("E49").Value = "  -0.59%  "
$ws.Range("E50").Value = "  -6.61%  "
$ws.Range("D51").Value = "119.20"
$ws.Range("E51").Value = "  -0.34%  "
